$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet.
#    We duplicate "2021-Q4" since it already has the exact same
#    layout/styles that the new quarter sheet needs, then we just
#    overwrite the data values.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$srcSheet.Copy($totalSheet)

# The freshly copied sheet is inserted immediately before "2021-Q4"'s
# original neighbor (i.e. right before "总计") and is auto-named
# "2021-Q4 (2)" by the engine - rename it.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Update the data row (header row, B2 "基金代码" and C2 "基金名称"
# already match the copied "2021-Q4" sheet's values, so only the
# figures that actually changed need to be written).
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "27.15"
$newSheet.Range("E2").Value = "86.02"
$newSheet.Range("F2").Value = "8.21"
$newSheet.Range("G2").Value = "2.2290"
$newSheet.Range("H2").Value = 1

# ------------------------------------------------------------------
# 2. Add a "2022-Q1" row at the top of the "总计" summary sheet and
#    shift the existing rows down, renumbering the index column (A).
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2.23

# Re-apply the index-column style (bordered/bold) to the new A2 cell,
# matching the style already used by the other index cells below it.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0

$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
